$d = $word.ActiveDocument

# Locate the two paragraphs we need by content, rather than a hard-coded
# index, so the script is resilient to any incidental paragraph shifts.
$mainPara = $null
$notePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs($i)
    $ppText = $pp.Range.Text
    if ($mainPara -eq $null -and $ppText.IndexOf('Data cleaning was an important') -ge 0) {
        $mainPara = $pp
    }
    if ($notePara -eq $null -and $ppText.IndexOf('Not sure what the real value of these parts are') -ge 0) {
        $notePara = $pp
    }
}

# ------------------------------------------------------------------
# 1) Rework the salary-anomaly discussion inside the "Data cleaning…"
#    paragraph: trim the long run, then rebuild the remainder as a
#    sequence of runs (two of them italicised: LowestSalary /
#    HighestSalary) describing the 0k -> 15k / 999k -> 300k change.
# ------------------------------------------------------------------

$pStart = $mainPara.Range.Start
$full = $mainPara.Range.Text

$startMarker = 'value of 999k (although'
$idxStart = $full.IndexOf($startMarker)
$endMarker = 'Moreover, salary values'
$idxEnd = $full.IndexOf($endMarker)

$segStart = $pStart + $idxStart
# +1 so the leading "M" of "Moreover" is swallowed too - it gets
# reinserted as the tail of the " M" run below, with the remaining
# "oreover, salary values…" text picking up right after it, matching
# the target run layout.
$segEnd = $pStart + $idxEnd + 1

$target = $d.Range($segStart, $segEnd)

$newText = 'value of 999k (although none had both a LowestSalary of 0k and HighestSalary of 999k). These salaries appeared to be anomalies as no reasonable salary pays no nothing, nor do so many pay exactly 999k. In addition, the SEEK website allows job searches from 30k to 200k+. As such, to keep closer to these amounts 0k LowestSalary rows were converted to 15k and 999k HighestSalary rows were converted to 300k. Accounting for these anomalies will help reduce inaccurate skewing of the data. M'

$target.Text = $newText

# Italicise the second occurrence of "LowestSalary" (the one followed by
# " rows were converted to 15k") to match the new emphasised term.
$scan = $d.Range($target.Start, $target.End)
[void]$scan.Find.Execute('LowestSalary rows were converted to 15k')
$scan.End = $scan.Start + 'LowestSalary'.Length
$scan.Italic = 1

# Italicise the second occurrence of "HighestSalary" (the one followed by
# " rows were converted to 300k") to match the new emphasised term.
$scan2 = $d.Range($target.Start, $target.End)
[void]$scan2.Find.Execute('HighestSalary rows were converted to 300k')
$scan2.End = $scan2.Start + 'HighestSalary'.Length
$scan2.Italic = 1

# ------------------------------------------------------------------
# 2) Remove the throw-away "# Not sure what the real value …" note
#    paragraph entirely (including its paragraph mark).
# ------------------------------------------------------------------

$noteRange = $d.Range($notePara.Range.Start, $notePara.Range.End)
$noteRange.Delete()

Write-Output 'done'
